# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
# for each changed row, matching the source commit's scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.926.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.07%  "

$ws.Range("D3").Value = "'3.031.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.81%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'593.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.66%  "

$ws.Range("D6").Value = "'153.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.66%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'3.028.03"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.77%  "

$ws.Range("D9").Value = "'0.516"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").Value = "'6.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +13.71%  "

$ws.Range("E11").Value = "  +4.21%  "

$ws.Range("E12").Value = "  +2.18%  "

$ws.Range("E13").Value = "  +3.53%  "

$ws.Range("D14").Value = "'35.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.90%  "

$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").Value = "'3.533.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.94%  "

$ws.Range("E17").Value = "  +2.88%  "

$ws.Range("D18").Value = "'62.904.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.07%  "

$ws.Range("D19").Value = "'3.031.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.86%  "

$ws.Range("D20").Value = "'452.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.33%  "

$ws.Range("E21").Value = "  +0.78%  "

$ws.Range("E22").Value = "  +2.69%  "

$ws.Range("E23").Value = "  +3.11%  "

$ws.Range("D24").Value = "'83.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.43%  "

$ws.Range("D25").Value = "'11.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.10%  "

$ws.Range("E26").Value = "  +5.81%  "

$ws.Range("E27").Value = "  +5.01%  "

$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("D29").Value = "'7.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.59%  "

$ws.Range("D30").Value = "'2.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.74%  "

$ws.Range("E31").Value = "  +0.79%  "

$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.18%  "

$ws.Range("E33").Value = "  +1.85%  "

$ws.Range("E34").Value = "  +1.88%  "

$ws.Range("D35").Value = "'0.0₃0858"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.78%  "

$ws.Range("E36").Value = "  +2.59%  "

$ws.Range("D37").Value = "'5.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.15%  "

$ws.Range("D38").Value = "'3.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.20%  "

$ws.Range("E39").Value = "  +7.44%  "

$ws.Range("D40").Value = "'2.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.33%  "

$ws.Range("D41").Value = "'50.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.43%  "

$ws.Range("D42").Value = "'9.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.25%  "

$ws.Range("D43").Value = "'0.308"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +15.92%  "

$ws.Range("D44").Value = "'44.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.47%  "

$ws.Range("D45").Value = "'390.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.44%  "

$ws.Range("E46").Value = "  +3.60%  "

$ws.Range("D47").Value = "'2.721.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.57%  "

$ws.Range("D48").Value = "'133.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.28%  "

$ws.Range("D50").Value = "'2.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.77%  "

$ws.Range("D51").Value = "'24.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.37%  "

